$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells are treated as text so numeric-looking
# strings (e.g. "177.51", "1.00") are preserved exactly, matching the
# source data which stores prices as plain text, not numbers.

# Updated cryptos list data (prices, 1h volume %, and two reordered rows)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.687.44"
$ws.Range("E2").Value = "  +2.59%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.532.83"
$ws.Range("E3").Value = "  +1.11%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.25"
$ws.Range("E5").Value = "  +2.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.51"
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.518"
$ws.Range("E8").Value = "  +0.61%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.530.87"
$ws.Range("E9").Value = "  +1.08%  "
$ws.Range("E10").Value = "  +12.90%  "
$ws.Range("E11").Value = "  -0.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.345"
$ws.Range("E12").Value = "  +1.22%  "
$ws.Range("E13").Value = "  +1.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000183"
$ws.Range("E14").Value = "  +6.13%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.982.39"
$ws.Range("E15").Value = "  +0.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.12"
$ws.Range("E16").Value = "  +1.61%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.526.59"
$ws.Range("E17").Value = "  +2.66%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.492.23"
$ws.Range("E18").Value = "  +0.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.67"
$ws.Range("E19").Value = "  +1.94%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "364.68"
$ws.Range("E20").Value = "  +3.35%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.09"
$ws.Range("E21").Value = "  +0.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.09"
$ws.Range("E22").Value = "  -0.49%  "
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.64"
$ws.Range("E24").Value = "  -0.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.26"
$ws.Range("E25").Value = "  -1.90%  "
$ws.Range("B26").Value = "SuiNetwork"
$ws.Range("C26").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.71"
$ws.Range("E26").Value = "  -2.49%  "
$ws.Range("B27").Value = "Aptos"
$ws.Range("C27").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.15"
$ws.Range("E27").Value = "  +0.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.657.49"
$ws.Range("E28").Value = "  +2.57%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  +0.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0927"
$ws.Range("E30").Value = "  +1.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "511.08"
$ws.Range("E31").Value = "  +1.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.76"
$ws.Range("E32").Value = "  -1.44%  "
$ws.Range("E33").Value = "  -1.23%  "
$ws.Range("E34").Value = "  +0.72%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("E36").Value = "  -1.80%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "161.56"
$ws.Range("E37").Value = "  -1.88%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.81"
$ws.Range("E38").Value = "  +2.18%  "
$ws.Range("E39").Value = "  +1.33%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.33"
$ws.Range("E40").Value = "  -1.10%  "
$ws.Range("E41").Value = "  +0.11%  "
$ws.Range("E42").Value = "  -0.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.87"
$ws.Range("E43").Value = "  -0.54%  "
$ws.Range("E44").Value = "  -2.28%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.43"
$ws.Range("E45").Value = "  -1.82%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "38.82"
$ws.Range("E46").Value = "  -0.44%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "151.10"
$ws.Range("E47").Value = "  +4.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.60"
$ws.Range("E48").Value = "  +1.73%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.519"
$ws.Range("E49").Value = "  -0.12%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0₆0253"
$ws.Range("E50").Value = "  -1.95%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0739"
$ws.Range("E51").Value = "  -0.47%  "
